# Updates cryptos list: refreshed prices/volumes, and swapped the
# WrappedEther/Chainlink row order (rows 15-16) to match the new ranking.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.951.82"
$ws.Range("D3").Value = "1.813.23"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'310.57"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").Value = "'0.4976"
$ws.Range("E7").Value = "  -3.00%  "
$ws.Range("D8").Value = "'0.3890"
$ws.Range("E8").Value = "  +2.88%  "
$ws.Range("D9").Value = "'0.09688"
$ws.Range("E9").Value = "  +24.66%  "
$ws.Range("E10").Value = "  +1.63%  "
$ws.Range("D11").Value = "'41.07"
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("D12").Value = "'6.423"
$ws.Range("E13").Value = "  +2.07%  "
$ws.Range("D14").Value = "'1.002"
$ws.Range("E14").Value = "  -0.02%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.811.28"
$ws.Range("E15").Value = "  +1.98%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'7.300"
$ws.Range("E16").Value = "  +2.02%  "
$ws.Range("E17").Value = "  +5.14%  "
$ws.Range("D18").Value = "'92.61"
$ws.Range("E18").Value = "  +1.22%  "
$ws.Range("D19").Value = "'0.06674"
$ws.Range("E19").Value = "  +2.20%  "
$ws.Range("E21").Value = "  +0.88%  "
$ws.Range("D22").Value = "'5.917"
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").Value = "28.001.99"
$ws.Range("E23").Value = "  +0.60%  "
$ws.Range("E24").Value = "  +1.41%  "
$ws.Range("D25").Value = "'2.242"
$ws.Range("E25").Value = "  +0.40%  "
$ws.Range("D26").Value = "'159.19"
$ws.Range("E26").Value = "  +0.29%  "
$ws.Range("D27").Value = "2.018.24"
$ws.Range("E27").Value = "  +1.82%  "
$ws.Range("D28").Value = "'20.54"
$ws.Range("E28").Value = "  +1.87%  "
$ws.Range("D29").Value = "'2.393"
$ws.Range("E29").Value = "  +2.03%  "
$ws.Range("E30").Value = "  +2.77%  "
$ws.Range("D31").Value = "'0.1060"
$ws.Range("E31").Value = "  -1.54%  "
$ws.Range("D32").Value = "'1.037"
$ws.Range("E32").Value = "  +1.23%  "
$ws.Range("D33").Value = "'5.567"
$ws.Range("E33").Value = "  +1.65%  "
$ws.Range("D34").Value = "'3.634"
$ws.Range("E34").Value = "  +0.51%  "
$ws.Range("D35").Value = "'0.06709"
$ws.Range("E35").Value = "  -4.71%  "
$ws.Range("D36").Value = "'9.005"
$ws.Range("E36").Value = "  +3.93%  "
$ws.Range("D37").Value = "'0.02328"
$ws.Range("E37").Value = "  +1.02%  "
$ws.Range("D38").Value = "'0.2135"
$ws.Range("E38").Value = "  +0.91%  "
$ws.Range("D39").Value = "'4.940"
$ws.Range("E39").Value = "  -1.47%  "
$ws.Range("E40").Value = "  -2.27%  "
$ws.Range("D41").Value = "'0.6186"
$ws.Range("E41").Value = "  +1.79%  "
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("E43").Value = "  -0.23%  "
$ws.Range("D44").Value = "'13.11"
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").Value = "'0.5871"
$ws.Range("E45").Value = "  -1.39%  "
$ws.Range("E46").Value = "  -0.69%  "
$ws.Range("D47").Value = "'1.280"
$ws.Range("E47").Value = "  -3.20%  "
$ws.Range("D48").Value = "'122.97"
$ws.Range("E48").Value = "  -3.10%  "
$ws.Range("D49").Value = "'1.937"
$ws.Range("E49").Value = "  +2.26%  "
$ws.Range("E50").Value = "  -2.56%  "
$ws.Range("D51").Value = "'0.06792"
$ws.Range("E51").Value = "  +1.37%  "
